# Applies the "Updated cryptos list" price/volume refresh described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells that are being rewritten to stay text-typed
# (several look like plain numbers, e.g. "331.11", and Excel's COM Value
# setter would otherwise auto-coerce numeric-looking strings to floats,
# losing the original text formatting). Each contiguous block is set in its
# own statement -- this engine only honours NumberFormat on the first area
# of a comma-joined multi-area Range.
$ws.Range("D2:D5").NumberFormat = "@"
$ws.Range("D7:D11").NumberFormat = "@"
$ws.Range("D13:D18").NumberFormat = "@"
$ws.Range("D20:D27").NumberFormat = "@"
$ws.Range("D29:D33").NumberFormat = "@"
$ws.Range("D35:D37").NumberFormat = "@"
$ws.Range("D39:D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '24.874.10'
$ws.Range("E2").Value = '  +1.64%  '

$ws.Range("D3").Value = '1.670.74'
$ws.Range("E3").Value = '  +0.84%  '

$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").Value = '331.11'
$ws.Range("E5").Value = '  +7.48%  '

$ws.Range("E6").Value = '  +0.06%  '

$ws.Range("D7").Value = '0.3648'
$ws.Range("E7").Value = '  +0.58%  '

$ws.Range("D8").Value = '46.81'
$ws.Range("E8").Value = '  -1.05%  '

$ws.Range("D9").Value = '0.3220'
$ws.Range("E9").Value = '  -1.48%  '

$ws.Range("D10").Value = '1.139'
$ws.Range("E10").Value = '  +1.32%  '

$ws.Range("D11").Value = '0.07046'
$ws.Range("E11").Value = '  +1.22%  '

$ws.Range("E12").Value = '  +0.22%  '

$ws.Range("D13").Value = '6.066'
$ws.Range("E13").Value = '  +2.40%  '

$ws.Range("D14").Value = '19.55'
$ws.Range("E14").Value = '  +1.22%  '

$ws.Range("D15").Value = '1.664.56'
$ws.Range("E15").Value = '  +0.67%  '

$ws.Range("D16").Value = '6.609'
$ws.Range("E16").Value = '  -0.06%  '

$ws.Range("D17").Value = '0.00001044'
$ws.Range("E17").Value = '  +0.28%  '

$ws.Range("D18").Value = '0.06534'
$ws.Range("E18").Value = '  +0.27%  '

$ws.Range("E19").Value = '  +0.17%  '

$ws.Range("D20").Value = '78.48'
$ws.Range("E20").Value = '  +2.78%  '

$ws.Range("D21").Value = '15.80'
$ws.Range("E21").Value = '  +0.58%  '

$ws.Range("D22").Value = '5.912'
$ws.Range("E22").Value = '  +0.08%  '

$ws.Range("D23").Value = '12.88'
$ws.Range("E23").Value = '  +2.42%  '

$ws.Range("D24").Value = '24.884.48'
$ws.Range("E24").Value = '  +1.88%  '

$ws.Range("D25").Value = '2.442'
$ws.Range("E25").Value = '  -0.71%  '

$ws.Range("D26").Value = '2.392'
$ws.Range("E26").Value = '  +3.50%  '

$ws.Range("D27").Value = '148.28'
$ws.Range("E27").Value = '  +1.40%  '

$ws.Range("E28").Value = '  +1.33%  '

$ws.Range("D29").Value = '1.849.10'
$ws.Range("E29").Value = '  +0.56%  '

$ws.Range("D30").Value = '125.42'
$ws.Range("E30").Value = '  +1.03%  '

$ws.Range("D31").Value = '1.170'
$ws.Range("E31").Value = '  -1.97%  '

$ws.Range("D32").Value = '4.074'
$ws.Range("E32").Value = '  +0.31%  '

$ws.Range("D33").Value = '5.722'
$ws.Range("E33").Value = '  +2.56%  '

$ws.Range("E34").Value = '  +1.10%  '

$ws.Range("D35").Value = '1.641'
$ws.Range("E35").Value = '  -2.54%  '

$ws.Range("D36").Value = '12.23'
$ws.Range("E36").Value = '  -0.89%  '

$ws.Range("D37").Value = '5.143'
$ws.Range("E37").Value = '  -1.21%  '

$ws.Range("E38").Value = '  +2.07%  '

$ws.Range("D39").Value = '0.06031'
$ws.Range("E39").Value = '  -0.31%  '

$ws.Range("D40").Value = '0.02229'
$ws.Range("E40").Value = '  +1.57%  '

$ws.Range("D41").Value = '0.2080'
$ws.Range("E41").Value = '  +1.55%  '

$ws.Range("D42").Value = '8.199'
$ws.Range("E42").Value = '  +0.33%  '

$ws.Range("D43").Value = '0.9997'
$ws.Range("E43").Value = '  -0.04%  '

$ws.Range("D44").Value = '0.5930'
$ws.Range("E44").Value = '  +1.14%  '

$ws.Range("D45").Value = '13.61'
$ws.Range("E45").Value = '  +7.49%  '

$ws.Range("D46").Value = '3.857'
$ws.Range("E46").Value = '  +3.35%  '

$ws.Range("D47").Value = '0.5698'
$ws.Range("E47").Value = '  +2.13%  '

$ws.Range("D48").Value = '124.64'
$ws.Range("E48").Value = '  +2.23%  '

$ws.Range("D49").Value = '1.955'
$ws.Range("E49").Value = '  +0.79%  '

$ws.Range("E50").Value = '  +1.34%  '

$ws.Range("D51").Value = '1.186'
$ws.Range("E51").Value = '  +3.12%  '
